# fix: fixed formatting when scrapping floating point numbers
#
# 1) Two "Razon social"/"Nombre Fantasia" names had a comma typo'd where a
#    period was intended.
# 2) The "Importe" column (H, rows 2-100) held Argentine-formatted amounts
#    stored as text (e.g. "1.140,00" = thousands "." + decimal ",").
#    They need to become plain decimal text (e.g. "1140.00": no thousands
#    separator, "." as the decimal mark) while STAYING text (not becoming
#    real numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Name fixes -----------------------------------------------------------
$ws.Range("E52").Value = "RICCOTTI. MARIANA EDITH"

$ws.Range("E68").Value = "GIMENEZ. ROBERTO ADRIAN"
$ws.Range("F68").Value = "GIMENEZ. ROBERTO ADRIAN"

# --- Importe column reformat (H2:H100) ------------------------------------
# Force text entry (so "1140.00" isn't auto-coerced into the number 1140),
# write the reformatted amount, then drop the temporary text number-format
# so the cell style stays the original default (General / unstyled).
for ($r = 2; $r -le 100; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $old = $cell.Text
    $new = $old.Replace(".", "").Replace(",", ".")
    $cell.NumberFormat = "@"
    $cell.Value = $new
    $cell.ClearFormats()
}
